$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = [double]"0.002027356019707161"
$ws.Range("I2").Value = [double]"0.002027356019707161"
$ws.Range("J2").Value = [double]"0.5891811476386277"
$ws.Range("K2").Value = [double]"0.5891811476386277"
$ws.Range("L2").Value = [double]"41.5403630465008"
$ws.Range("M2").Value = "[15.171865799070773, 67.90886029393083]"
$ws.Range("N2").Value = [double]"0.002719205218807197"
$ws.Range("O2").Value = [double]"0.002719205218807197"
$ws.Range("P2").Value = [double]"1.415131825941349"
$ws.Range("Q2").Value = "[0.6352369529781159, 2.1950266989045817]"
$ws.Range("R2").Value = [double]"0.0006705025592619318"
$ws.Range("S2").Value = [double]"0.0006705025592619318"
$ws.Range("T2").Value = [double]"61.47517743107858"
$ws.Range("U2").Value = "[46.37019881854863, 76.58015604360853]"
$ws.Range("V2").Value = [double]"1.783686531808826e-10"
$ws.Range("W2").Value = [double]"1.783686531808826e-10"
$ws.Range("X2").Value = [double]"17.04504504504504"
$ws.Range("Y2").Value = [double]"14.31431431431431"
$ws.Range("Z2").Value = [double]"19.77577577577578"

# Row 3
$ws.Range("H3").Value = [double]"0.001895382316740468"
$ws.Range("I3").Value = [double]"0.001895382316740468"
$ws.Range("J3").Value = [double]"0.1696113817097465"
$ws.Range("K3").Value = [double]"0.1696113817097465"
$ws.Range("L3").Value = [double]"35.94388283543005"
$ws.Range("M3").Value = "[11.702255356496686, 60.18551031436342]"
$ws.Range("N3").Value = [double]"0.004555558110825642"
$ws.Range("O3").Value = [double]"0.004555558110825642"
$ws.Range("P3").Value = [double]"2.06923720326535"
$ws.Range("Q3").Value = "[1.3396581285578106, 2.7988162779728887]"
$ws.Range("R3").Value = [double]"8.328616902808506e-07"
$ws.Range("S3").Value = [double]"8.328616902808506e-07"
$ws.Range("T3").Value = [double]"52.65278541856107"
$ws.Range("U3").Value = "[39.24034674334139, 66.06522409378076]"
$ws.Range("V3").Value = [double]"4.713904822040149e-10"
$ws.Range("W3").Value = [double]"4.713904822040149e-10"
$ws.Range("X3").Value = [double]"14.75475475475475"
$ws.Range("Y3").Value = [double]"12.2002002002002"
$ws.Range("Z3").Value = [double]"17.30930930930931"

# Row 4
$ws.Range("H4").Value = [double]"4.202640178818662e-05"
$ws.Range("I4").Value = [double]"4.202640178818662e-05"
$ws.Range("J4").Value = [double]"0.7411492416803547"
$ws.Range("K4").Value = [double]"0.7411492416803547"
$ws.Range("L4").Value = [double]"44.12716220211635"
$ws.Range("M4").Value = "[19.7682797902812, 68.4860446139515]"
$ws.Range("N4").Value = [double]"0.0006826490933802276"
$ws.Range("O4").Value = [double]"0.0006826490933802276"
$ws.Range("P4").Value = [double]"2.547237286694427"
$ws.Range("Q4").Value = "[2.0063424554457336, 3.08813211794312]"
$ws.Range("R4").Value = [double]"2.652988939644274e-12"
$ws.Range("S4").Value = [double]"2.652988939644274e-12"
$ws.Range("T4").Value = [double]"53.09647059125567"
$ws.Range("U4").Value = "[40.290962445255644, 65.90197873725569]"
$ws.Range("V4").Value = [double]"1.068289900985064e-10"
$ws.Range("W4").Value = [double]"1.068289900985064e-10"
$ws.Range("X4").Value = [double]"13.08108108108108"
$ws.Range("Y4").Value = [double]"11.18718718718719"
$ws.Range("Z4").Value = [double]"14.97497497497497"

# Row 5
$ws.Range("F5").Value = [double]"22.55000000000009"
$ws.Range("H5").Value = [double]"7.689429639334833e-05"
$ws.Range("I5").Value = [double]"7.689429639334833e-05"
$ws.Range("L5").Value = [double]"45.32243924209806"
$ws.Range("M5").Value = "[19.472745766352148, 71.17213271784397]"
$ws.Range("N5").Value = [double]"0.0009677944313333242"
$ws.Range("O5").Value = [double]"0.0009677944313333242"
$ws.Range("P5").Value = [double]"2.471763589310888"
$ws.Range("Q5").Value = "[1.9182898084982725, 3.025237370123504]"
$ws.Range("R5").Value = [double]"1.288702478063897e-11"
$ws.Range("S5").Value = [double]"1.288702478063897e-11"
$ws.Range("T5").Value = [double]"82.58243122813727"
$ws.Range("U5").Value = "[69.1124648476638, 96.05239760861073]"
$ws.Range("V5").Value = [double]"4.440892098500626e-16"
$ws.Range("W5").Value = [double]"4.440892098500626e-16"
$ws.Range("X5").Value = [double]"13.67897897897903"
$ws.Range("Y5").Value = [double]"11.69259259259264"
$ws.Range("Z5").Value = [double]"15.66536536536542"

# Row 6
$ws.Range("F6").Value = [double]"22.55000000000009"
$ws.Range("H6").Value = [double]"2.534067289339248e-08"
$ws.Range("I6").Value = [double]"2.534067289339248e-08"
$ws.Range("L6").Value = [double]"59.63004872921704"
$ws.Range("M6").Value = "[37.987079211667094, 81.27301824676698]"
$ws.Range("N6").Value = [double]"1.450230371347061e-06"
$ws.Range("O6").Value = [double]"1.450230371347061e-06"
$ws.Range("P6").Value = [double]"2.836553126664658"
$ws.Range("Q6").Value = "[2.434026740619119, 3.2390795127101963]"
$ws.Range("R6").Value = [double]"0"
$ws.Range("S6").Value = [double]"0"
$ws.Range("T6").Value = [double]"67.71392014100674"
$ws.Range("U6").Value = "[55.72642266583294, 79.70141761618054]"
$ws.Range("V6").Value = [double]"7.771561172376096e-15"
$ws.Range("W6").Value = [double]"7.771561172376096e-15"
$ws.Range("X6").Value = [double]"12.36976976976982"
$ws.Range("Y6").Value = [double]"10.92512512512517"
$ws.Range("Z6").Value = [double]"13.81441441441447"

# Row 7
$ws.Range("F7").Value = [double]"22.55000000000009"
$ws.Range("H7").Value = [double]"0.0007496098217150138"
$ws.Range("I7").Value = [double]"0.0007496098217150138"
$ws.Range("L7").Value = [double]"39.78549070378454"
$ws.Range("M7").Value = "[17.580920960047884, 61.990060447521195]"
$ws.Range("N7").Value = [double]"0.0007689693720362989"
$ws.Range("O7").Value = [double]"0.0007689693720362989"
$ws.Range("P7").Value = [double]"-3.119579491852928"
$ws.Range("Q7").Value = "[-3.886895415252237, -2.352263568453619]"
$ws.Range("R7").Value = [double]"1.835793739246583e-10"
$ws.Range("S7").Value = [double]"1.835793739246583e-10"
$ws.Range("T7").Value = [double]"49.58951665567889"
$ws.Range("U7").Value = "[35.46472484281361, 63.71430846854416]"
$ws.Range("V7").Value = [double]"7.993209427681336e-09"
$ws.Range("W7").Value = [double]"7.993209427681336e-09"
$ws.Range("X7").Value = [double]"11.19599599599604"
$ws.Range("Y7").Value = [double]"8.442142142142174"
$ws.Range("Z7").Value = [double]"13.9498498498499"
